$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3157.6
$ws.Range("J17").Value = 3290.2856
$ws.Range("L17").Value = 9870.856800000001
$ws.Range("N17").Value = -10206.8568

$ws.Range("H29").Value = 4828.2856
$ws.Range("J29").Value = 5199.8335
$ws.Range("L29").Value = 15599.5005
$ws.Range("N29").Value = -16161.5005

$ws.Range("H32").Value = 789.1
$ws.Range("I32").Value = 749.6
$ws.Range("K32").Value = 749.6
$ws.Range("M32").Value = -423.6

$ws.Range("H33").Value = 230.46666
$ws.Range("I33").Value = 62.46154
$ws.Range("K33").Value = 62.46154
$ws.Range("M33").Value = 166.53846

$ws.Range("H70").Value = 19475.5
$ws.Range("I70").Value = 1499
$ws.Range("J70").Value = 25467.666
$ws.Range("K70").Value = 4497
$ws.Range("L70").Value = 76402.99800000001
$ws.Range("M70").Value = -4227
$ws.Range("N70").Value = -76942.99800000001

$ws.Range("H73").Value = 19475.5
$ws.Range("I73").Value = 1499
$ws.Range("J73").Value = 25467.666
$ws.Range("K73").Value = 4497
$ws.Range("L73").Value = 76402.99800000001
$ws.Range("M73").Value = -3561
$ws.Range("N73").Value = -78274.99800000001

$ws.Range("H116").Value = 4611.25
$ws.Range("I116").Value = 1999.25
$ws.Range("K116").Value = 1999.25
$ws.Range("M116").Value = 1442.75

$ws.Range("H129").Value = 3130.8462
$ws.Range("I129").Value = 1464.1666
$ws.Range("K129").Value = 4392.4998
$ws.Range("M129").Value = 607.5002000000004

$ws.Range("H138").Value = 13697.462
$ws.Range("I138").Value = 13098.2
$ws.Range("J138").Value = 14072
$ws.Range("K138").Value = 39294.60000000001
$ws.Range("L138").Value = 42216
$ws.Range("M138").Value = -34154.60000000001
$ws.Range("N138").Value = -52496

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 2996
$ws.Range("I28").Value = 2996
$ws.Range("K28").Value = 2996
$ws.Range("M28").Value = -2804

$ws.Range("H30").Value = 4222
$ws.Range("I30").Value = 4222
$ws.Range("K30").Value = 4222
$ws.Range("M30").Value = -4072

$ws.Range("H32").Value = 14726.621
$ws.Range("I32").Value = 6395.2573
$ws.Range("J32").Value = 24133
$ws.Range("K32").Value = 6395.2573
$ws.Range("L32").Value = 24133
$ws.Range("M32").Value = -6108.2573
$ws.Range("N32").Value = -24707

$ws.Range("H99").Value = 2996
$ws.Range("I99").Value = 2996
$ws.Range("K99").Value = 2996
$ws.Range("M99").Value = -1

$ws.Range("H102").Value = 1518.75
$ws.Range("I102").Value = 1295.7273
$ws.Range("J102").Value = 2009.4
$ws.Range("K102").Value = 1295.7273
$ws.Range("L102").Value = 2009.4
$ws.Range("M102").Value = 326.2727
$ws.Range("N102").Value = -5253.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H97").Value = 13156.333
$ws.Range("I97").Value = 13156.333
$ws.Range("K97").Value = 13156.333
$ws.Range("M97").Value = -12165.333

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 4422.5
$ws.Range("I105").Value = 2673.5
$ws.Range("K105").Value = 2673.5
$ws.Range("M105").Value = -926.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 1028.4
$ws.Range("I75").Value = 659
$ws.Range("J75").Value = 1120.75
$ws.Range("K75").Value = 1977
$ws.Range("L75").Value = 3362.25
$ws.Range("M75").Value = -979
$ws.Range("N75").Value = -5358.25

$ws.Range("H78").Value = 1028.4
$ws.Range("I78").Value = 659
$ws.Range("J78").Value = 1120.75
$ws.Range("K78").Value = 5931
$ws.Range("L78").Value = 10086.75
$ws.Range("M78").Value = -939
$ws.Range("N78").Value = -20070.75

$ws.Range("H113").Value = 2167.8
$ws.Range("J113").Value = 1864.3334
$ws.Range("L113").Value = 5593.0002
$ws.Range("N113").Value = -9933.0002

$ws.Range("H131").Value = 1479.7333
$ws.Range("I131").Value = 622.8
$ws.Range("J131").Value = 1908.2
$ws.Range("K131").Value = 1868.4
$ws.Range("L131").Value = 5724.6
$ws.Range("M131").Value = 3171.6
$ws.Range("N131").Value = -15804.6

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 8408.223
$ws.Range("I70").Value = 6891.6665
$ws.Range("J70").Value = 9166.5
$ws.Range("K70").Value = 6891.6665
$ws.Range("L70").Value = 9166.5
$ws.Range("M70").Value = -6621.6665
$ws.Range("N70").Value = -9706.5

$ws.Range("H73").Value = 8408.223
$ws.Range("I73").Value = 6891.6665
$ws.Range("J73").Value = 9166.5
$ws.Range("K73").Value = 6891.6665
$ws.Range("L73").Value = 9166.5
$ws.Range("M73").Value = -5955.6665
$ws.Range("N73").Value = -11038.5

$ws.Range("H80").Value = 3838.2222
$ws.Range("I80").Value = 2655.75
$ws.Range("K80").Value = 2655.75
$ws.Range("M80").Value = -1657.75

$ws.Range("H83").Value = 3838.2222
$ws.Range("I83").Value = 2655.75
$ws.Range("K83").Value = 13278.75
$ws.Range("M83").Value = -8286.75

$ws.Range("H113").Value = 3156.8333
$ws.Range("I113").Value = 2670
$ws.Range("J113").Value = 3254.2
$ws.Range("K113").Value = 2670
$ws.Range("L113").Value = 3254.2
$ws.Range("M113").Value = -500
$ws.Range("N113").Value = -7594.2

$ws.Range("H122").Value = 613827.4
$ws.Range("I122").Value = 85492.086
$ws.Range("J122").Value = 1670498
$ws.Range("K122").Value = 256476.258
$ws.Range("L122").Value = 5011494
$ws.Range("M122").Value = -254026.258
$ws.Range("N122").Value = -5016394

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H104").Value = 30073.6
$ws.Range("J104").Value = 30073.6
$ws.Range("L104").Value = 30073.6
$ws.Range("N104").Value = -37061.6

$ws.Range("H136").Value = 4080.625
$ws.Range("I136").Value = 3949.2856
$ws.Range("K136").Value = 11847.8568
$ws.Range("M136").Value = -9297.856800000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 2500375
$ws.Range("J2").Value = 2500375
$ws.Range("L2").Value = 2500375
$ws.Range("N2").Value = -2500599

$ws.Range("H45").Value = 21633
$ws.Range("I45").Value = 25449.5
$ws.Range("J45").Value = 14000
$ws.Range("K45").Value = 25449.5
$ws.Range("L45").Value = 14000
$ws.Range("M45").Value = -24958.5
$ws.Range("N45").Value = -14982

$ws.Range("H96").Value = 941.4
$ws.Range("I96").Value = 967.6667
$ws.Range("J96").Value = 902
$ws.Range("K96").Value = 967.6667
$ws.Range("L96").Value = 902
$ws.Range("M96").Value = 405.3333
$ws.Range("N96").Value = -3648

$ws.Range("H122").Value = 1401.4286
$ws.Range("I122").Value = 1370.5264
$ws.Range("J122").Value = 1695
$ws.Range("K122").Value = 4111.5792
$ws.Range("L122").Value = 5085
$ws.Range("M122").Value = -1661.5792
